$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T45")

# New values for rows 2-11, columns B-G
$data = @{
    2  = @(0.2517957102541956, 1.386019378421211, 7.589302934471694, 2.754868950507754, 2.80499345858655, 23)
    3  = @(0.5243254742780024, 1.822206725996549, 13.98149530182482, 3.739183774813004, 3.789363045952605, 22)
    4  = @(-0.8019598128076217, 1.305564707247611, 6.943070400165851, 2.634970664004791, 2.571950116496777, 21)
    5  = @(-0.02187051981424342, 0.4503981798184682, 0.5039077891406407, 0.709864627334424, 0.7279599970962077, 20)
    6  = @(0.0766489273505621, 0.6798153727021657, 1.082404327345544, 1.040386623974734, 1.065990830592977, 19)
    7  = @(-0.1907918648194715, 0.5798188246083275, 0.8145864660070502, 0.9025444399070055, 0.9077227079139136, 18)
    8  = @(-0.1371052455939829, 0.4571923072893577, 0.4943098393230912, 0.7030717170553024, 0.7107963775852487, 17)
    9  = @(0.1145279375814358, 0.4251814381380697, 0.3204394163604665, 0.5660736845680662, 0.5725477730274827, 16)
    10 = @(-0.01358233361013406, 0.3321431173215987, 0.2545945103064831, 0.5045735925576001, 0.5220940294212741, 15)
    11 = @(-0.02420063924481428, 0.3457558232543146, 0.2682231761632309, 0.517902670550395, 0.5384613548423456, 13)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
